$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.077415413507486
$ws.Range("D2").Value = 1.082194288545775
$ws.Range("E2").Value = 1.090580024769537
$ws.Range("F2").Value = 1.096518449284907
$ws.Range("I2").Value = 1.064255357879657
$ws.Range("J2").Value = 1.082309987651964
$ws.Range("K2").Value = 1.084863140559048
$ws.Range("L2").Value = 1.093227215020482
$ws.Range("M2").Value = 1.099150515618444
$ws.Range("N2").Value = 1.03074174736267

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.078672277550381
$ws.Range("D3").Value = 1.083243992898239
$ws.Range("E3").Value = 1.091799918735069
$ws.Range("F3").Value = 1.097751408378859
$ws.Range("I3").Value = 1.064717765019582
$ws.Range("J3").Value = 1.083225614113197
$ws.Range("K3").Value = 1.085730872144224
$ws.Range("L3").Value = 1.094266297436399
$ws.Range("M3").Value = 1.100203729181428
$ws.Range("N3").Value = 1.031073092788312

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.079485084481155
$ws.Range("D4").Value = 1.083922790667796
$ws.Range("E4").Value = 1.092589134044327
$ws.Range("F4").Value = 1.09854913720466
$ws.Range("I4").Value = 1.065015452722452
$ws.Range("J4").Value = 1.083817057006806
$ws.Range("K4").Value = 1.086291313351568
$ws.Range("L4").Value = 1.094937938810367
$ws.Range("M4").Value = 1.100884574784167
$ws.Range("N4").Value = 1.031286595030291

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.079826679091176
$ws.Range("D5").Value = 1.084208055652723
$ws.Range("E5").Value = 1.092920889261736
$ws.Range("F5").Value = 1.098884485720113
$ws.Range("I5").Value = 1.065140237616334
$ws.Range("J5").Value = 1.084065455116721
$ws.Range("K5").Value = 1.086526675659646
$ws.Range("L5").Value = 1.095220127703526
$ws.Range("M5").Value = 1.101170647574911
$ws.Range("N5").Value = 1.031376136333549

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.079884027989711
$ws.Range("D6").Value = 1.084255946988185
$ws.Range("E6").Value = 1.092976590670394
$ws.Range("F6").Value = 1.098940791292066
$ws.Range("I6").Value = 1.065161168291488
$ws.Range("J6").Value = 1.084107147967499
$ws.Range("K6").Value = 1.086566179578793
$ws.Range("L6").Value = 1.095267498596979
$ws.Range("M6").Value = 1.101218671437808
$ws.Range("N6").Value = 1.031391158117475

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.079489649311689
$ws.Range("D7").Value = 1.083926602789671
$ws.Range("E7").Value = 1.092593567092612
$ws.Range("F7").Value = 1.098553618210949
$ws.Range("I7").Value = 1.065017121529434
$ws.Range("J7").Value = 1.083820377072931
$ws.Range("K7").Value = 1.08629445924348
$ws.Range("L7").Value = 1.094941710094403
$ws.Range("M7").Value = 1.100888397908329
$ws.Range("N7").Value = 1.031287792329868

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.077840275125638
$ws.Range("D8").Value = 1.082549131291417
$ws.Range("E8").Value = 1.090992323388841
$ws.Range("F8").Value = 1.096935150986745
$ws.Range("I8").Value = 1.064411946014102
$ws.Range("J8").Value = 1.082619641962358
$ws.Range("K8").Value = 1.08515661063781
$ws.Range("L8").Value = 1.093578526857254
$ws.Range("M8").Value = 1.099506590613177
$ws.Range("N8").Value = 1.030853913848129

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.074930162694742
$ws.Range("D9").Value = 1.080118475660398
$ws.Range("E9").Value = 1.08816958066938
$ws.Range("F9").Value = 1.094082515276396
$ws.Range("I9").Value = 1.063333857004134
$ws.Range("J9").Value = 1.080495835117668
$ws.Range("K9").Value = 1.083143545431339
$ws.Range("L9").Value = 1.091170866971592
$ws.Range("M9").Value = 1.097066576807512
$ws.Range("N9").Value = 1.030082445110749

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072987428879777
$ws.Range("D10").Value = 1.078495665290456
$ws.Range("E10").Value = 1.086286849636805
$ws.Range("F10").Value = 1.092180164661146
$ws.Range("I10").Value = 1.062607203858741
$ws.Range("J10").Value = 1.079074499623924
$ws.Range("K10").Value = 1.081796000127834
$ws.Range("L10").Value = 1.089561903088406
$ws.Range("M10").Value = 1.095436357279949
$ws.Range("N10").Value = 1.029563446719603

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.072145534978924
$ws.Range("D11").Value = 1.077792383640973
$ws.Range("E11").Value = 1.085471362235801
$ws.Range("F11").Value = 1.091356255452097
$ws.Range("I11").Value = 1.062290659758109
$ws.Range("J11").Value = 1.078457723715381
$ws.Range("K11").Value = 1.081211169735723
$ws.Range("L11").Value = 1.088864261435169
$ws.Range("M11").Value = 1.094729586433005
$ws.Range("N11").Value = 1.029337595444086

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.071832712965872
$ws.Range("D12").Value = 1.077531062115857
$ws.Range("E12").Value = 1.085168413299344
$ws.Range("F12").Value = 1.091050189370454
$ws.Range("I12").Value = 1.062172794592208
$ws.Range("J12").Value = 1.078228423804122
$ws.Range("K12").Value = 1.080993735103559
$ws.Range("L12").Value = 1.088604981111399
$ws.Range("M12").Value = 1.094466926362722
$ws.Range("N12").Value = 1.029253535002602

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.071899819161274
$ws.Range("D13").Value = 1.077587120641159
$ws.Range("E13").Value = 1.085233398755238
$ws.Range("F13").Value = 1.091115842962253
$ws.Range("I13").Value = 1.062198090061288
$ws.Range("J13").Value = 1.078277618586154
$ws.Range("K13").Value = 1.081040384771589
$ws.Range("L13").Value = 1.088660604229193
$ws.Range("M13").Value = 1.094523273938412
$ws.Range("N13").Value = 1.029271573921899

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.072119679177946
$ws.Range("D14").Value = 1.07777078459937
$ws.Range("E14").Value = 1.085446321226094
$ws.Range("F14").Value = 1.091330956551796
$ws.Range("I14").Value = 1.062280922841045
$ws.Range("J14").Value = 1.07843877384486
$ws.Range("K14").Value = 1.081193200666095
$ws.Range("L14").Value = 1.088842832218099
$ws.Range("M14").Value = 1.094707877615526
$ws.Range("N14").Value = 1.029330650438022

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.072255128158415
$ws.Range("D15").Value = 1.077883933848723
$ws.Range("E15").Value = 1.0855775043267
$ws.Range("F15").Value = 1.091463491143324
$ws.Range("I15").Value = 1.06233192083325
$ws.Range("J15").Value = 1.078538040098369
$ws.Range("K15").Value = 1.08128732867459
$ws.Range("L15").Value = 1.088955089588285
$ws.Range("M15").Value = 1.094821600216216
$ws.Range("N15").Value = 1.02936702698424

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.073043288047021
$ws.Range("D16").Value = 1.078542327077398
$ws.Range("E16").Value = 1.086340965302695
$ws.Range("F16").Value = 1.092234840808086
$ws.Range("I16").Value = 1.062628171732618
$ws.Range("J16").Value = 1.079115404826805
$ws.Range("K16").Value = 1.081834785101079
$ws.Range("L16").Value = 1.0896081830478
$ws.Range("M16").Value = 1.095483244688902
$ws.Range("N16").Value = 1.029578412033633

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.073537496201937
$ws.Range("D17").Value = 1.078955159189867
$ws.Range("E17").Value = 1.086819794727288
$ws.Range("F17").Value = 1.092718638326112
$ws.Range("I17").Value = 1.062813492850889
$ws.Range("J17").Value = 1.079477213217205
$ws.Range("K17").Value = 1.082177831314161
$ws.Range("L17").Value = 1.090017595321175
$ws.Range("M17").Value = 1.095898040790591
$ws.Range("N17").Value = 1.029710707510859

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.073825694333144
$ws.Range("D18").Value = 1.079195900006487
$ws.Range("E18").Value = 1.087099063655317
$ws.Range("F18").Value = 1.093000812158285
$ws.Range("I18").Value = 1.062921404439673
$ws.Range("J18").Value = 1.079688121930156
$ws.Range("K18").Value = 1.082377795761434
$ws.Range("L18").Value = 1.090256307112192
$ws.Range("M18").Value = 1.096139900077412
$ws.Range("N18").Value = 1.029787765060681

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073923951503376
$ws.Range("D19").Value = 1.079277976771347
$ws.Range("E19").Value = 1.087194283122946
$ws.Range("F19").Value = 1.093097023326354
$ws.Range("I19").Value = 1.062958168489554
$ws.Range("J19").Value = 1.079760014709646
$ws.Range("K19").Value = 1.082445956720118
$ws.Range("L19").Value = 1.090337686224941
$ws.Range("M19").Value = 1.096222353625233
$ws.Range("N19").Value = 1.029814021369018

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.073484479105837
$ws.Range("D20").Value = 1.078910872141213
$ws.Range("E20").Value = 1.086768423382649
$ws.Range("F20").Value = 1.092666733217526
$ws.Range("I20").Value = 1.062793628596498
$ws.Range("J20").Value = 1.079438407858449
$ws.Range("K20").Value = 1.082141039020468
$ws.Range("L20").Value = 1.089973678738129
$ws.Range("M20").Value = 1.095853545828805
$ws.Range("N20").Value = 1.029696524651688

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.072054938812197
$ws.Range("D21").Value = 1.077716702682024
$ws.Range("E21").Value = 1.085383622007871
$ws.Range("F21").Value = 1.09126761179736
$ws.Range("I21").Value = 1.062256538574305
$ws.Range("J21").Value = 1.078391323229319
$ws.Range("K21").Value = 1.081148205803303
$ws.Range("L21").Value = 1.088789174644063
$ws.Range("M21").Value = 1.094653520149949
$ws.Range("N21").Value = 1.02931325855362

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.071155520149552
$ws.Range("D22").Value = 1.076965351020455
$ws.Range("E22").Value = 1.084512706153337
$ws.Range("F22").Value = 1.090387756536729
$ws.Range("I22").Value = 1.061917189944373
$ws.Range("J22").Value = 1.077731810413327
$ws.Range("K22").Value = 1.080522799011265
$ws.Range("L22").Value = 1.088043588876986
$ws.Range("M22").Value = 1.093898240406123
$ws.Range("N22").Value = 1.029071304522858

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.071632377943637
$ws.Range("D23").Value = 1.077363707705974
$ws.Range("E23").Value = 1.0849744184189
$ws.Range("F23").Value = 1.090854201657063
$ws.Range("I23").Value = 1.062097242767344
$ws.Range("J23").Value = 1.078081542234343
$ws.Range("K23").Value = 1.08085445089136
$ws.Range("L23").Value = 1.088438918547103
$ws.Range("M23").Value = 1.09429870286034
$ws.Range("N23").Value = 1.02919966201197

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.073508435450626
$ws.Range("D24").Value = 1.078930883728855
$ws.Range("E24").Value = 1.086791635956732
$ws.Range("F24").Value = 1.092690186954632
$ws.Range("I24").Value = 1.062802604963697
$ws.Range("J24").Value = 1.079455942725192
$ws.Range("K24").Value = 1.082157664271643
$ws.Range("L24").Value = 1.089993523034229
$ws.Range("M24").Value = 1.09587365144712
$ws.Range("N24").Value = 1.029702933609485

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.075682952802415
$ws.Range("D25").Value = 1.080747269047652
$ws.Range("E25").Value = 1.088899476738393
$ws.Range("F25").Value = 1.094820085052761
$ws.Range("I25").Value = 1.063613961234858
$ws.Range("J25").Value = 1.081045844861111
$ws.Range("K25").Value = 1.083664933242496
$ws.Range("L25").Value = 1.0917939757355
$ws.Range("M25").Value = 1.097697994673671
$ws.Range("N25").Value = 1.030282711721707
